$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A58").NumberFormat = "@"
$ws.Range("A58").Value = "2025/10/04"
$ws.Range("A58").Style = "Normal"
$ws.Range("B58").Value = "土"
$ws.Range("C58").Value = 8
$ws.Range("D58").Value = 5
